$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($row, $col, $val)
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue 2 4 "51.881.73"
Set-TextValue 2 5 "  -1.04%  "

Set-TextValue 3 4 "2.925.79"
Set-TextValue 3 5 "  +0.02%  "

Set-TextValue 4 5 "  -0.02%  "

Set-TextValue 5 4 "360.52"
Set-TextValue 5 5 "  +2.25%  "

Set-TextValue 6 4 "110.34"
Set-TextValue 6 5 "  -2.48%  "

Set-TextValue 7 5 "  +0.98%  "

Set-TextValue 8 5 "  +0.01%  "

Set-TextValue 9 5 "  +0.39%  "

Set-TextValue 10 4 "39.37"
Set-TextValue 10 5 "  -2.58%  "

Set-TextValue 11 4 "0.0879"
Set-TextValue 11 5 "  +1.69%  "

Set-TextValue 12 5 "  +0.87%  "

Set-TextValue 13 4 "19.65"
Set-TextValue 13 5 "  -2.60%  "

Set-TextValue 14 4 "7.93"
Set-TextValue 14 5 "  +0.47%  "

Set-TextValue 15 4 "3.391.32"
Set-TextValue 15 5 "  +0.14%  "

Set-TextValue 16 4 "2.932.99"
Set-TextValue 16 5 "  -0.23%  "

Set-TextValue 17 4 "0.990"
Set-TextValue 17 5 "  -0.28%  "

Set-TextValue 18 4 "51.937.78"
Set-TextValue 18 5 "  -1.05%  "

Set-TextValue 19 5 "  +0.76%  "

Set-TextValue 20 4 "7.61"
Set-TextValue 20 5 "  -1.25%  "

Set-TextValue 21 5 "  -2.97%  "

Set-TextValue 22 4 "0.0₃0984"
Set-TextValue 22 5 "  -0.06%  "

Set-TextValue 23 5 "  -0.09%  "

Set-TextValue 24 5 "  -0.56%  "

Set-TextValue 25 5 "  +1.00%  "

Set-TextValue 26 5 "  +13.11%  "

Set-TextValue 27 4 "27.04"
Set-TextValue 27 5 "  -0.42%  "

Set-TextValue 28 4 "7.64"
Set-TextValue 28 5 "  +16.78%  "

Set-TextValue 29 4 "0.999"
Set-TextValue 29 5 "  -0.08%  "

Set-TextValue 30 5 "  +14.00%  "

Set-TextValue 31 5 "  -0.44%  "

Set-TextValue 32 4 "38.29"
Set-TextValue 32 5 "  +0.54%  "

Set-TextValue 33 5 "  +1.70%  "

Set-TextValue 34 4 "6.09"
Set-TextValue 34 5 "  -1.93%  "

Set-TextValue 35 4 "52.13"
Set-TextValue 35 5 "  -1.91%  "

Set-TextValue 36 5 "  -1.95%  "

Set-TextValue 37 5 "  +0.01%  "

Set-TextValue 38 5 "  -2.82%  "

Set-TextValue 39 4 "18.42"
Set-TextValue 39 5 "  -2.68%  "

Set-TextValue 40 4 "2.01"
Set-TextValue 40 5 "  -3.82%  "

Set-TextValue 41 4 "2.73"
Set-TextValue 41 5 "  -0.70%  "

Set-TextValue 42 5 "  +2.15%  "

Set-TextValue 43 4 "23.18"
Set-TextValue 43 5 "  -5.04%  "

Set-TextValue 44 4 "119.24"
Set-TextValue 44 5 "  -2.80%  "

Set-TextValue 45 5 "  -1.49%  "

Set-TextValue 46 5 "  -2.15%  "

Set-TextValue 47 2 "Maker"
Set-TextValue 47 3 "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue 47 4 "2.133.34"
Set-TextValue 47 5 "  -4.05%  "

Set-TextValue 48 2 "TheGraph"
Set-TextValue 48 3 "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue 48 4 "0.251"
Set-TextValue 48 5 "  -4.92%  "

Set-TextValue 49 2 "BEAM"
Set-TextValue 49 3 "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
Set-TextValue 49 4 "0.0335"
Set-TextValue 49 5 "  -2.00%  "

Set-TextValue 50 2 "FraxShare"
Set-TextValue 50 3 "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue 50 4 "9.19"
Set-TextValue 50 5 "  +0.11%  "

Set-TextValue 51 2 "SEI"
Set-TextValue 51 3 "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
Set-TextValue 51 4 "0.913"
Set-TextValue 51 5 "  -5.40%  "

